$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain number (e.g. "1.00", "28.20") must be
# forced to Text format first, otherwise Excel will store them as numeric
# values and normalize/round the display text (losing trailing zeros, etc).
$textCells = @("D5", "D6", "D7", "D9", "D10", "D11", "D12", "D13", "D14", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D28", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D41", "D42", "D43", "D44", "D45", "D47", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated coin data values
$ws.Range("D2").Value = "64.584.82"
$ws.Range("E2").Value = "  +5.58%  "
$ws.Range("D3").Value = "3.094.83"
$ws.Range("E3").Value = "  +4.09%  "
$ws.Range("E4").Value = "  +0.35%  "
$ws.Range("D5").Value = "558.17"
$ws.Range("E5").Value = "  +2.44%  "
$ws.Range("D6").Value = "144.35"
$ws.Range("E6").Value = "  +11.51%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("D8").Value = "3.075.96"
$ws.Range("E8").Value = "  +3.44%  "
$ws.Range("D9").Value = "0.499"
$ws.Range("E9").Value = "  +2.24%  "
$ws.Range("D10").Value = "7.05"
$ws.Range("E10").Value = "  +18.34%  "
$ws.Range("D11").Value = "0.152"
$ws.Range("E11").Value = "  +6.88%  "
$ws.Range("D12").Value = "0.460"
$ws.Range("E12").Value = "  +4.89%  "
$ws.Range("D13").Value = "35.34"
$ws.Range("E13").Value = "  +5.85%  "
$ws.Range("D14").Value = "0.0000226"
$ws.Range("E14").Value = "  +4.88%  "
$ws.Range("D15").Value = "3.596.81"
$ws.Range("E15").Value = "  +4.28%  "
$ws.Range("D16").Value = "64.767.49"
$ws.Range("E16").Value = "  +5.77%  "
$ws.Range("D17").Value = "3.106.00"
$ws.Range("E17").Value = "  +4.72%  "
$ws.Range("E18").Value = "  -0.46%  "
$ws.Range("D19").Value = "6.72"
$ws.Range("E19").Value = "  +2.86%  "
$ws.Range("D20").Value = "484.90"
$ws.Range("E20").Value = "  +2.27%  "
$ws.Range("D21").Value = "13.72"
$ws.Range("E21").Value = "  +5.97%  "
$ws.Range("D22").Value = "0.672"
$ws.Range("E22").Value = "  +2.43%  "
$ws.Range("D23").Value = "7.53"
$ws.Range("E23").Value = "  +9.18%  "
$ws.Range("D24").Value = "13.26"
$ws.Range("E24").Value = "  +11.81%  "
$ws.Range("D25").Value = "80.75"
$ws.Range("E25").Value = "  +1.81%  "
$ws.Range("E26").Value = "  +0.50%  "
$ws.Range("E27").Value = "  +5.72%  "
$ws.Range("D28").Value = "7.97"
$ws.Range("E28").Value = "  +5.58%  "
$ws.Range("E29").Value = "  +9.79%  "
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.33%  "
$ws.Range("D31").Value = "26.08"
$ws.Range("E31").Value = "  +3.34%  "
$ws.Range("D32").Value = "1.15"
$ws.Range("E32").Value = "  +3.31%  "
$ws.Range("D33").Value = "2.45"
$ws.Range("E33").Value = "  +8.28%  "
$ws.Range("D34").Value = "5.75"
$ws.Range("E34").Value = "  +6.85%  "
$ws.Range("D35").Value = "55.19"
$ws.Range("E35").Value = "  +1.32%  "
$ws.Range("D36").Value = "6.11"
$ws.Range("E36").Value = "  +5.44%  "
$ws.Range("D37").Value = "465.67"
$ws.Range("E37").Value = "  +4.91%  "
$ws.Range("D38").Value = "0.0408"
$ws.Range("E38").Value = "  +9.15%  "
$ws.Range("D39").Value = "0.0825"
$ws.Range("E39").Value = "  +5.63%  "
$ws.Range("D40").Value = "3.038.24"
$ws.Range("E40").Value = "  -1.61%  "
$ws.Range("D41").Value = "0.117"
$ws.Range("E41").Value = "  +2.37%  "
$ws.Range("D42").Value = "8.29"
$ws.Range("E42").Value = "  +3.83%  "
$ws.Range("D43").Value = "2.70"
$ws.Range("E43").Value = "  +18.51%  "
$ws.Range("D44").Value = "28.20"
$ws.Range("E44").Value = "  +12.71%  "
$ws.Range("D45").Value = "0.260"
$ws.Range("E45").Value = "  +9.46%  "
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").Value = "2.08"
$ws.Range("E47").Value = "  +9.00%  "
$ws.Range("E48").Value = "  +5.02%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "118.18"
$ws.Range("E49").Value = "  +3.88%  "
$ws.Range("B50").Value = "PEPE"
$ws.Range("C50").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D50").Value = "0.0₃0517"
$ws.Range("E50").Value = "  +9.46%  "
$ws.Range("D51").Value = "2.07"
$ws.Range("E51").Value = "  +5.38%  "
